$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New title row (row 1) — did not exist before (dimension started at A2),
#    so simply writing into row 1 does not shift any existing rows.
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").Merge()
$ws.Range("A1").Value = "TABLA1"
$ws.Range("E1:F1").Merge()
$ws.Range("E1").Value = "TABLA2"

# ---------------------------------------------------------------------------
# 2. Row 2 caption text change ("REGRESIÓN" -> "DATOS REGRESIÓN MULTIPLE").
#    E2 ("VENTANA DE PARZEN") keeps its text.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "DATOS REGRESIÓN MULTIPLE"

# ---------------------------------------------------------------------------
# 3. Fill in the right-hand "TABLA2" data (F4:F7) with the ECM results and
#    make sure the E4:E7 smoothing-window labels carry through (already
#    present from the source file, values unchanged).
# ---------------------------------------------------------------------------
$ws.Range("F4").Value = 28163.0528
$ws.Range("F5").Value = 27880.9913
$ws.Range("F6").Value = 30281.0351
$ws.Range("F7").Value = 32611.2016

# ---------------------------------------------------------------------------
# 4. Fonts: whole used block becomes Times New Roman 8pt; the unused tail of
#    the right-hand table (E8:F13) is intentionally left at the old default
#    (Calibri 11) because that sub-table only has 4 data rows. Column D is
#    never touched (it has no cells in either version), so every range below
#    is split at D to avoid materialising empty D cells.
# ---------------------------------------------------------------------------
foreach ($rng in @($ws.Range("A1:C7"), $ws.Range("E1:F7"), $ws.Range("A8:C13"))) {
    $rng.Font.Name = "Times New Roman"
    $rng.Font.Size = 8
}

# Bold for the three header rows, regular weight for the data rows.
foreach ($rng in @($ws.Range("A1:C3"), $ws.Range("E1:F3"))) {
    $rng.Font.Bold = $true
}
foreach ($rng in @($ws.Range("A4:C7"), $ws.Range("E4:F7"), $ws.Range("A8:C13"))) {
    $rng.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 5. Alignment: headers centred horizontally (already true for most cells,
#    but row 1 / row 2 right-hand blank cells need it too).
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("E1:F1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 6. Borders.
#    Row 3 (column headers): top+bottom thin, no left/right.
#    Row 2 (captions): bottom thin only.
#    Rows 4-13 (data): no borders at all.
# ---------------------------------------------------------------------------
foreach ($rng in @($ws.Range("A3:C3"), $ws.Range("E3:F3"))) {
    $rng.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft
    $rng.Borders.Item(10).LineStyle = -4142  # xlEdgeRight
    $rng.Borders.Item(11).LineStyle = -4142  # xlInsideVertical
    $rng.Borders.Item(8).LineStyle = 1       # xlEdgeTop
    $rng.Borders.Item(9).LineStyle = 1       # xlEdgeBottom
}

foreach ($rng in @($ws.Range("A2:C2"), $ws.Range("E2:F2"))) {
    $rng.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft
    $rng.Borders.Item(10).LineStyle = -4142  # xlEdgeRight
    $rng.Borders.Item(11).LineStyle = -4142  # xlInsideVertical
    $rng.Borders.Item(8).LineStyle = -4142   # xlEdgeTop
    $rng.Borders.Item(9).LineStyle = 1       # xlEdgeBottom
}

$ws.Range("A4:C13").Borders.LineStyle = -4142
$ws.Range("E4:F13").Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# 7. Number formats: C column already carried 0.0000; extend it to F (the
#    new ECM results column) for the header + all data rows.
# ---------------------------------------------------------------------------
$ws.Range("F3:F13").NumberFormat = "0.0000"

# ---------------------------------------------------------------------------
# 8. Column F width/style + selection, matching the refreshed layout.
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 11.42578125

$ws.Range("F12").Select()
